# Progress.xlsx update — fill in missing lab scores on sheet "БИВТ-22-17"
# (the cells that previously held a stray "pass" text label get a real
# numeric score), mark E25 as "failed 3,7", and leave the cursor on the
# sheet/cell the author ended on.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- numeric lab scores that replace the placeholder "pass" text ---
$ws1.Range("F2").Value  = 5
$ws1.Range("D6").Value  = 5
$ws1.Range("D7").Value  = 5
$ws1.Range("E7").Value  = 5
$ws1.Range("E9").Value  = 5
$ws1.Range("D10").Value = 5
$ws1.Range("E10").Value = 5
$ws1.Range("E16").Value = 5
$ws1.Range("E19").Value = 5
$ws1.Range("F20").Value = 5
$ws1.Range("D21").Value = 5
$ws1.Range("H23").Value = 5
$ws1.Range("F28").Value = 5
$ws1.Range("F31").Value = 4
$ws1.Range("C32").Value = 5

# --- E25 becomes an explicit failure note instead of the generic "pass" text ---
$ws1.Range("E25").Value = "failed 3,7"

# --- the author finished on sheet1, with E26 selected ---
$ws1.Activate()
$ws1.Range("E26").Select()
